$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 988440.0600000001
$ws.Range("E2").Value = 642249.28
$ws.Range("H2").Value = 28696.86
$ws.Range("I2").Value = 28880
$ws.Range("J2").Value = 28132.53
$ws.Range("K2").Value = 978027.0600000001

$ws.Range("D3").Value = 779343.73
$ws.Range("E3").Value = 618612.09
$ws.Range("G3").Value = 23651.42
$ws.Range("H3").Value = 21889.07
$ws.Range("I3").Value = 23652.02
$ws.Range("J3").Value = 20368.05
$ws.Range("K3").Value = 774208.73
$ws.Range("L3").Value = -203818.33
$ws.Range("M3").Value = -39.69
$ws.Range("N3").Value = -23637.19
$ws.Range("Q3").Value = -203818.54

$ws.Range("D4").Value = 875395.59
$ws.Range("E4").Value = 644601.73
$ws.Range("G4").Value = 26007.39
$ws.Range("H4").Value = 25054.73
$ws.Range("I4").Value = 26007.39
$ws.Range("J4").Value = 23807.98
$ws.Range("K4").Value = 868128.59
$ws.Range("L4").Value = -109898.47
$ws.Range("M4").Value = -15.12
$ws.Range("N4").Value = 2352.45
$ws.Range("Q4").Value = -109898.47

$ws.Range("D5").Value = 914190.96
$ws.Range("E5").Value = 648365.89
$ws.Range("G5").Value = 27027.39
$ws.Range("H5").Value = 26276.06
$ws.Range("I5").Value = 27027.39
$ws.Range("J5").Value = 25152.01
$ws.Range("K5").Value = 906442.96
$ws.Range("L5").Value = -71584.10000000001
$ws.Range("M5").Value = -9.24
$ws.Range("N5").Value = 6116.62
$ws.Range("Q5").Value = -71584.10000000001

$ws.Range("D6").Value = 927221.49
$ws.Range("E6").Value = 651093.13
$ws.Range("G6").Value = 27327.39
$ws.Range("H6").Value = 26721.94
$ws.Range("I6").Value = 27327.39
$ws.Range("J6").Value = 25690.13
$ws.Range("K6").Value = 918940.49
$ws.Range("L6").Value = -59086.57
$ws.Range("M6").Value = -7.14
$ws.Range("N6").Value = 8843.85
$ws.Range("Q6").Value = -59086.57

$ws.Range("D7").Value = 952888.71
$ws.Range("E7").Value = 652032.4300000001
$ws.Range("G7").Value = 27986.92
$ws.Range("H7").Value = 27542.8
$ws.Range("I7").Value = 27986.92
$ws.Range("J7").Value = 26643.45
$ws.Range("K7").Value = 944074.71
$ws.Range("L7").Value = -33952.35
$ws.Range("M7").Value = -3.85
$ws.Range("N7").Value = 9783.15
$ws.Range("Q7").Value = -33952.35

$ws.Range("D8").Value = 961365.6800000001
$ws.Range("E8").Value = 650206.1
$ws.Range("G8").Value = 28170.36
$ws.Range("H8").Value = 27842.67
$ws.Range("I8").Value = 28170.36
$ws.Range("J8").Value = 27059.48
$ws.Range("K8").Value = 952018.6800000001
$ws.Range("L8").Value = -26008.39
$ws.Range("M8").Value = -2.78
$ws.Range("N8").Value = 7956.82
$ws.Range("Q8").Value = -26008.39

$ws.Range("D9").Value = 983133.02
$ws.Range("E9").Value = 647245.53
$ws.Range("G9").Value = 28767.58
$ws.Range("H9").Value = 28507.09
$ws.Range("I9").Value = 28767.58
$ws.Range("J9").Value = 27820.36
$ws.Range("K9").Value = 973253.02
$ws.Range("L9").Value = -4774.04
$ws.Range("M9").Value = -0.48
$ws.Range("N9").Value = 4996.25
$ws.Range("Q9").Value = -4774.04

$ws.Range("D10").Value = 992164.86
$ws.Range("E10").Value = 635670.78
$ws.Range("G10").Value = 28952.64
$ws.Range("H10").Value = 28835.29
$ws.Range("I10").Value = 28952.64
$ws.Range("J10").Value = 28387.21
$ws.Range("K10").Value = 981181.86
$ws.Range("L10").Value = 3154.79
$ws.Range("M10").Value = 0.29
$ws.Range("N10").Value = -6578.5
$ws.Range("Q10").Value = 3154.79

$ws.Range("D11").Value = 1010257.07
$ws.Range("E11").Value = 629035.08
$ws.Range("G11").Value = 29460.77
$ws.Range("H11").Value = 29377.69
$ws.Range("I11").Value = 29460.77
$ws.Range("J11").Value = 29011.37
$ws.Range("K11").Value = 998700.0699999999
$ws.Range("L11").Value = 20673.01
$ws.Range("M11").Value = 1.79
$ws.Range("N11").Value = -13214.2
$ws.Range("Q11").Value = 20673.01

$ws.Range("D12").Value = 1010991.04
$ws.Range("E12").Value = 619465.74
$ws.Range("G12").Value = 29462.45
$ws.Range("H12").Value = 29415.55
$ws.Range("I12").Value = 29462.45
$ws.Range("J12").Value = 29145.14
$ws.Range("K12").Value = 998864.04
$ws.Range("L12").Value = 20836.97
$ws.Range("M12").Value = 1.72
$ws.Range("N12").Value = -22783.53
$ws.Range("Q12").Value = 20836.97

$ws.Range("D13").Value = 1022006.61
$ws.Range("E13").Value = 605753.4
$ws.Range("G13").Value = 29773.32
$ws.Range("H13").Value = 29744.54
$ws.Range("I13").Value = 29773.32
$ws.Range("J13").Value = 29559.52
$ws.Range("K13").Value = 1009305.61
$ws.Range("L13").Value = 31278.54
$ws.Range("M13").Value = 2.46
$ws.Range("N13").Value = -36495.88
$ws.Range("Q13").Value = 31278.54

$ws.Range("D14").Value = 1021304.89
$ws.Range("E14").Value = 594748.38
$ws.Range("G14").Value = 29747.27
$ws.Range("H14").Value = 29728.82
$ws.Range("I14").Value = 29747.27
$ws.Range("J14").Value = 29597.87
$ws.Range("K14").Value = 1008033.89
$ws.Range("L14").Value = 30006.82
$ws.Range("M14").Value = 2.26
$ws.Range("N14").Value = -47500.9
$ws.Range("Q14").Value = 30006.82

$ws.Range("D15").Value = 1039806.01
$ws.Range("E15").Value = 578218.29
$ws.Range("G15").Value = 30283.59
$ws.Range("H15").Value = 30269.5
$ws.Range("I15").Value = 30283.59
$ws.Range("J15").Value = 30163.36
$ws.Range("K15").Value = 1026015.01
$ws.Range("L15").Value = 47987.94
$ws.Range("M15").Value = 3.48
$ws.Range("N15").Value = -64030.99
$ws.Range("Q15").Value = 47987.94

$ws.Range("D16").Value = 1045947.06
$ws.Range("E16").Value = 514296.93
$ws.Range("G16").Value = 30461.55
$ws.Range("H16").Value = 30449.02
$ws.Range("I16").Value = 30461.55
$ws.Range("J16").Value = 30354.6
$ws.Range("K16").Value = 1029868.06
$ws.Range("L16").Value = 51841
$ws.Range("M16").Value = 3.22
$ws.Range("N16").Value = -127952.35
$ws.Range("Q16").Value = 51841

$ws.Range("D17").Value = 915582.22
$ws.Range("E17").Value = 565238.42
$ws.Range("H17").Value = 24792.95
$ws.Range("I17").Value = 27993.65
$ws.Range("J17").Value = 23871.44
$ws.Range("K17").Value = 905169.22
$ws.Range("L17").Value = -72857.84
$ws.Range("M17").Value = -7
$ws.Range("N17").Value = -77010.86
$ws.Range("Q17").Value = -77010.86
